$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Type" column (M) records whether a well contained deuterium-labelled
# ants (L) or unlabelled ants (U). Wells 7, 4 and 1 were mislabelled as
# unlabelled ("U") when they actually held labelled ants, so flip the Type
# column for those rows to "L".
$rows = @(11,12,13,14,15,16,17,18,19,20,21,22,29,30,31,32,33,34,47,48,49,50,51,52)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 13).Value = "L"
}

# Leave the sheet scrolled to the bottom of the data with the cell just
# past the last row selected, matching where the editor ended up after
# reviewing/making the change.
$ws.Range("A26").Select()
$excel.ActiveWindow.ScrollRow = 26
$ws.Range("M53").Select()
